$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "69.807.03"
Set-TextValue "E2" "  +2.49%  "

Set-TextValue "D3" "3.413.25"
Set-TextValue "E3" "  +2.05%  "

Set-TextValue "E4" "  -0.06%  "

Set-TextValue "D5" "587.70"
Set-TextValue "E5" "  +0.45%  "

Set-TextValue "D6" "181.87"
Set-TextValue "E6" "  +3.27%  "

Set-TextValue "E7" "  +1.05%  "

Set-TextValue "E8" "  +0.11%  "

Set-TextValue "E9" "  +11.25%  "

Set-TextValue "E10" "  +2.06%  "

Set-TextValue "D11" "48.68"
Set-TextValue "E11" "  +2.11%  "

Set-TextValue "E12" "  +5.29%  "

Set-TextValue "D13" "689.35"
Set-TextValue "E13" "  -1.25%  "

Set-TextValue "D14" "8.71"
Set-TextValue "E14" "  +3.50%  "

Set-TextValue "D15" "3.962.82"
Set-TextValue "E15" "  +1.77%  "

Set-TextValue "D16" "69.817.86"
Set-TextValue "E16" "  +2.40%  "

Set-TextValue "B17" "WrappedEther"
Set-TextValue "C17" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D17" "3.434.31"
Set-TextValue "E17" "  +2.64%  "

Set-TextValue "B18" "TRON"
Set-TextValue "C18" "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue "D18" "0.121"
Set-TextValue "E18" "  +1.61%  "

Set-TextValue "D19" "17.83"
Set-TextValue "E19" "  +1.91%  "

Set-TextValue "E20" "  +2.29%  "

Set-TextValue "D21" "0.917"
Set-TextValue "E21" "  +2.21%  "

Set-TextValue "D22" "17.38"
Set-TextValue "E22" "  +2.01%  "

Set-TextValue "E23" "  -0.47%  "

Set-TextValue "D24" "103.46"
Set-TextValue "E24" "  +2.05%  "

Set-TextValue "E25" "  +0.90%  "

Set-TextValue "D26" "2.73"
Set-TextValue "E26" "  +1.09%  "

Set-TextValue "D27" "9.79"
Set-TextValue "E27" "  +3.40%  "

Set-TextValue "E28" "  +2.68%  "

Set-TextValue "E29" "  +3.62%  "

Set-TextValue "E30" "  -0.51%  "

Set-TextValue "B31" "Cosmos"
Set-TextValue "C31" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D31" "11.19"
Set-TextValue "E31" "  +1.41%  "

Set-TextValue "B32" "Bittensor"
Set-TextValue "C32" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D32" "559.03"
Set-TextValue "E32" "  -3.03%  "

Set-TextValue "B33" "dogwifhat"
Set-TextValue "C33" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D33" "3.67"
Set-TextValue "E33" "  +9.33%  "

Set-TextValue "E34" "  +1.57%  "

Set-TextValue "D35" "58.72"
Set-TextValue "E35" "  +3.18%  "

Set-TextValue "E36" "  -0.05%  "

Set-TextValue "D37" "3.673.65"
Set-TextValue "E37" "  -2.40%  "

Set-TextValue "D38" "0.142"
Set-TextValue "E38" "  +5.82%  "

Set-TextValue "D39" "36.14"
Set-TextValue "E39" "  +2.01%  "

Set-TextValue "D40" "0.0₃0744"
Set-TextValue "E40" "  +9.46%  "

Set-TextValue "E41" "  +4.15%  "

Set-TextValue "D42" "2.71"
Set-TextValue "E42" "  +3.25%  "

Set-TextValue "B43" "VeChain"
Set-TextValue "C43" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D43" "0.0433"
Set-TextValue "E43" "  +6.18%  "

Set-TextValue "B44" "TheGraph"
Set-TextValue "C44" "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue "D44" "0.341"
Set-TextValue "E44" "  +1.73%  "

Set-TextValue "D45" "3.36"
Set-TextValue "E45" "  +2.00%  "

Set-TextValue "E46" "  +1.65%  "

Set-TextValue "E47" "  +0.76%  "

Set-TextValue "E48" "  +4.24%  "

Set-TextValue "D49" "1.00"
Set-TextValue "E49" "  -0.27%  "

Set-TextValue "D50" "130.83"
Set-TextValue "E50" "  +0.33%  "

Set-TextValue "D51" "2.71"
Set-TextValue "E51" "  +0.60%  "

